$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to Text format so purely-numeric-looking values
# (e.g. "579.09", "7.00") are preserved exactly as strings, not coerced to numbers.
foreach ($addr in @("D5","D6","D10","D15","D20","D21","D22","D26","D27","D29","D30","D31","D36","D37","D38","D40","D41","D42","D43","D47","D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.596.25'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '3.443.17'
$ws.Range("E3").Value = '  +2.49%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '579.09'
$ws.Range("E5").Value = '  +1.24%  '
$ws.Range("D6").Value = '146.90'
$ws.Range("E6").Value = '  +7.50%  '
$ws.Range("D7").Value = '3.444.07'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D10").Value = '7.64'
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").Value = '4.033.07'
$ws.Range("E13").Value = '  +2.50%  '
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").Value = '27.69'
$ws.Range("E15").Value = '  +7.05%  '
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").Value = '3.446.31'
$ws.Range("E17").Value = '  +2.52%  '
$ws.Range("D18").Value = '61.714.87'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("E19").Value = '  +7.73%  '
$ws.Range("D20").Value = '14.06'
$ws.Range("E20").Value = '  +0.93%  '
$ws.Range("D21").Value = '9.47'
$ws.Range("E21").Value = '  +1.85%  '
$ws.Range("D22").Value = '386.87'
$ws.Range("E22").Value = '  +2.93%  '
$ws.Range("E23").Value = '  +2.18%  '
$ws.Range("D24").Value = '3.590.34'
$ws.Range("E24").Value = '  +2.33%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = '5.77'
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '72.21'
$ws.Range("E27").Value = '  +1.69%  '
$ws.Range("E28").Value = '  -1.03%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '0.175'
$ws.Range("E29").Value = '  +7.03%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '7.81'
$ws.Range("E30").Value = '  +4.06%  '
$ws.Range("D31").Value = '1.58'
$ws.Range("E31").Value = '  -11.40%  '
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").Value = '24.20'
$ws.Range("E36").Value = '  +2.79%  '
$ws.Range("D37").Value = '5.23'
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("D38").Value = '7.00'
$ws.Range("E38").Value = '  +2.90%  '
$ws.Range("E39").Value = '  +1.96%  '
$ws.Range("D40").Value = '166.15'
$ws.Range("E40").Value = '  +0.82%  '
$ws.Range("D41").Value = '0.0788'
$ws.Range("E41").Value = '  +3.01%  '
$ws.Range("D42").Value = '25.91'
$ws.Range("E42").Value = '  +8.26%  '
$ws.Range("D43").Value = '0.789'
$ws.Range("E43").Value = '  +2.37%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("E45").Value = '  +2.79%  '
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").Value = '42.04'
$ws.Range("E47").Value = '  +1.60%  '
$ws.Range("D48").Value = '2.629.45'
$ws.Range("E48").Value = '  +10.98%  '
$ws.Range("E49").Value = '  -2.63%  '
$ws.Range("D50").Value = '23.84'
$ws.Range("E50").Value = '  +4.35%  '
$ws.Range("E51").Value = '  +0.25%  '
